$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix F22: "-" -> "PSU" (Salinity unit corrected)
$ws.Range("F22").Value = "PSU"

# Re-add UnitRegexp column (G), replacing the old UnitSynonyms column
$ws.Range("G1").Value = "UnitRegexp"
$ws.Range("G4").Value = "%s"
$ws.Range("G5").Value = "mass %"
$ws.Range("G6").Value = "mass %"
$ws.Range("G7").Value = "mass %"
$ws.Range("G8").Value = "mass %"
$ws.Range("G9").Value = "mass %"
$ws.Range("G10").Value = "mass %"
$ws.Range("G11").Value = "vol %"
$ws.Range("G12").Value = "mass %"
$ws.Range("G13").Value = "(µ|u)mol(\sL(⁻¹|-1)|/L)"
$ws.Range("G14").Value = "(µ|u)mol(\sL(⁻¹|-1)|/L)"
$ws.Range("G15").Value = "(µ|u)mol(\sL(⁻¹|-1)|/L)"
$ws.Range("G16").Value = "(µ|u)mol(\sL(⁻¹|-1)|/L)"
$ws.Range("G17").Value = "(µ|u)mol(\sL(⁻¹|-1)|/L)"
$ws.Range("G18").Value = "V"
$ws.Range("G19").Value = "(µ|u)mol(\sL(⁻¹|-1)|/L)"
$ws.Range("G20").Value = "(°|deg\s?)C"
$ws.Range("G21").Value = "mS(\scm(⁻¹|-1)|/cm)"
$ws.Range("G22").Value = "(PSU|psu)"
$ws.Range("G23").Value = "(L|l)ux"
$ws.Range("G24").Value = "(C|c)ounts"
$ws.Range("G25").Value = "(C|c)ounts"
$ws.Range("G26").Value = "d(|eci)b(|ar)"
$ws.Range("G27").Value = "(NTU|ntu)"
$ws.Range("G28").Value = "kg(\sm(⁻³|-3)|(/m|/m³)"
$ws.Range("G29").Value = "m(\ss(⁻¹|-1)|/s)"
$ws.Range("G30").Value = "m(\ss(⁻¹|-1)|/s)"
$ws.Range("G31").Value = "m(\ss(⁻¹|-1)|/s)"
$ws.Range("G32").Value = "mmol(\sL(⁻¹|-1)|/L)"
$ws.Range("G33").Value = "mmol(\sL(⁻¹|-1)|/L)"
$ws.Range("G34").Value = "mmol(\sL(⁻¹|-1)|/L)"

# Set column G width (Excel's ColumnWidth property has an offset vs. the
# serialized OOXML "width" attribute; 28.1666... round-trips to width="29")
$ws.Columns("G").ColumnWidth = 28.166666666666668

# Update view: scroll/selection state
$ws.Range("G28").Select()
